$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Range("M1").Value = "description"
$ws.Range("M1").Select()
